$d = $word.ActiveDocument

# Locate the body paragraph "Hyper parameter tuning in deep learning"
# (NOT the "Heading 1" styled title paragraph just above it, which has
# the same words but Title Case).
$target = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*Hyper parameter tuning in deep learning*" -and $p.Style.NameLocal -ne "Heading 1") {
        $target = $p
        break
    }
}

# The paragraph right after it is the stray empty paragraph ("<w:p/>")
# that should be removed as part of this edit.
$afterP = $target.Next()

$replaceRange = $d.Range($target.Range.Start, $afterP.Range.End)

$newXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:after="0"/></w:pPr><w:r><w:t>ann.add(Dense(4,input_dim = 6, activation="relu"))</w:t></w:r><w:r><w:tab/></w:r><w:r><w:tab/><w:t>use max no of hidden layer</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:after="0"/></w:pPr><w:r><w:t>ann.add(Dense(3, activation="relu"))</w:t></w:r><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:tab/><w:t>no of nodes in parameter shapes</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:after="0"/></w:pPr><w:r><w:t>ann.add(Dense(2, activation="relu"))</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:after="0"/></w:pPr><w:r><w:t>ann.add(Dense(1, activation="sigmoid"))</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>batch_size=100</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>epochs = 50</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>optimizer='adam'</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
'@

$replaceRange.InsertXML($newXml)
